# Daily attendance processing - 2026-01-20 22:01:19
# Reorders the "Recorded By" (column G) entries for each attendance row:
# the last name/address in the comma-separated list is moved to the front
# (a right-rotation of the list), except when the value is exactly
# "admin@admin.com, System", which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text -eq "admin@admin.com, System") {
        continue
    }

    if ($text.Contains(",")) {
        $parts = $text.Split(",")
        $trimmedParts = @()
        foreach ($p in $parts) {
            $trimmedParts += $p.Trim()
        }

        $count = $trimmedParts.Count
        $lastItem = $trimmedParts[$count - 1]

        $newOrder = @($lastItem)
        for ($k = 0; $k -le $count - 2; $k++) {
            $newOrder += $trimmedParts[$k]
        }

        $newValue = [string]::Join(", ", $newOrder)
        $cell.Value = $newValue
    }
}
